$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add column M, mirroring column L's formatting, then overwrite values ---
$ws.Range("L2:L15").Copy()
$ws.Range("M2:M15").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("M3").Value = 2020

$ws.Range("M4").Value = 94.1
$ws.Range("M5").Value = 99.6
$ws.Range("M6").Value = 91
$ws.Range("M7").Value = 86.886172668979881
$ws.Range("M8").Value = 86.955790296225956
$ws.Range("M9").Value = 96.29195112324031
$ws.Range("M10").Value = 97.849780305474511
$ws.Range("M11").Value = 90.676703333930902
$ws.Range("M12").Value = 99.675929342188979
$ws.Range("M13").Value = 100
$ws.Range("M14").Value = 100
$ws.Range("M15").Value = 100

# Row 4's new figure is emphasised (bold) relative to the rest of the column
$ws.Range("M4").Font.Bold = $true

# --- View state: scroll one column right, select G15 ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("G15").Select()
